$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -0.05188361557894466, 0.3732691606683305, 0.2921056387580342, 0.5404679812514653, 0.5435468534395178, 49),
    @(3, 0.1398461492215967, 0.5201991840348531, 0.5256151831718032, 0.7249932297420461, 0.7189057200126604, 48),
    @(4, 0.4203416940300701, 0.7015093178813666, 0.9741072080335899, 0.9869686965824144, 0.902433627878138, 48),
    @(5, 0.4178943810801925, 0.7896986526005529, 1.360587987782973, 1.166442449408874, 1.10078782933836, 47),
    @(6, 0.2532779346943542, 0.5959877371359621, 0.6944367481990631, 0.8333287155733103, 0.8026788987968445, 46),
    @(7, 0.2847813693415759, 0.5189089493205867, 0.4842484560603196, 0.6958796275652275, 0.6444878855629137, 34),
    @(8, 0.2226667514831955, 0.5151721287749901, 0.504580525438979, 0.7103383175916804, 0.6849954704767146, 33),
    @(9, 0.1530436297492116, 0.5048367829134828, 0.4592703115339658, 0.6776948513409009, 0.681839049572768, 16),
    @(10, 0.1712605454146585, 0.3442148425105949, 0.1719451111591538, 0.4146626474125127, 0.3980716528792657, 10)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
